# "added more selenium stuff" - add test_suite + OpenAccountTest sheets,
# extend AddCustomerTest with more sample rows.

$wb = $excel.ActiveWorkbook
$wsAddCustomer = $wb.Worksheets.Item(1)

# Create the two new sheets in the right tab order:
# AddCustomerTest, test_suite, OpenAccountTest
$wsTestSuite = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsAddCustomer)
$wsTestSuite.Name = "test_suite"

$wsOpenAccount = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsTestSuite)
$wsOpenAccount.Name = "OpenAccountTest"

# ---- OpenAccountTest sheet data ----
$wsOpenAccount.Range("A1").Value = "customer"
$wsOpenAccount.Range("B1").Value = "currency"
$wsOpenAccount.Range("A2").Value = "Nick Hero"
$wsOpenAccount.Range("B2").Value = "Rupee"

# ---- AddCustomerTest sheet: extra sample rows ----
$wsAddCustomer.Range("A3").Value = "Nick1"
$wsAddCustomer.Range("A4").Value = "Nick2"
$wsAddCustomer.Range("A5").Value = "Nick3"
$wsAddCustomer.Range("B3").Value = "Hero1"
$wsAddCustomer.Range("B4").Value = "Hero2"
$wsAddCustomer.Range("B5").Value = "Hero3"
$wsAddCustomer.Range("C3").Value = 223098
$wsAddCustomer.Range("C4").Value = 223098
$wsAddCustomer.Range("C5").Value = 223098
$wsAddCustomer.Range("D3").Value = "Customer added successfully"
$wsAddCustomer.Range("D4").Value = "Customer added successfully"
$wsAddCustomer.Range("D5").Value = "Customer added successfully"

# ---- test_suite sheet data ----
$wsTestSuite.Range("A1").Value = "TCID"
$wsTestSuite.Range("B1").Value = "Runmode"
$wsTestSuite.Range("A2").Value = "BankManagerLoginTest"
$wsTestSuite.Range("A3").Value = "AddCustomerTest"
$wsTestSuite.Range("A4").Value = "OpenAccountTest"
$wsTestSuite.Range("B2").Value = "Y"
$wsTestSuite.Range("B3").Value = "Y"
$wsTestSuite.Range("B4").Value = "Y"

# Font for the TCID column values (Menlo 9, black) - build the style once
# on A2, then replicate the exact same cell formatting onto A3:A4 via
# PasteSpecial so all three share a single cellXf/style entry.
$wsTestSuite.Range("A2").Font.Name = "Menlo"
$wsTestSuite.Range("A2").Font.Size = 9
$wsTestSuite.Range("A2").Font.Color = 0
$wsTestSuite.Range("A2").Copy()
$wsTestSuite.Range("A3:A4").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# ---- Column widths ----
# (ColumnWidth is in characters; Excel pads it by ~5/6 of a character when
# writing the XML `width` attribute, so back the character value off to
# land on the intended stored width of 34.)
$wsTestSuite.Columns.Item(1).ColumnWidth = 33.166666666666664

# ---- Selections per sheet ----
$wsAddCustomer.Range("D15").Select()
$wsTestSuite.Range("C10").Select()
$wsOpenAccount.Range("C9").Select()

# test_suite is the active tab when the workbook is opened
$wsTestSuite.Activate()
